# Applies the "Added Tsai, organized files, and Added basic html tags" commit
# to MouseDatasets.xlsx: shortens / re-organizes the dataset "Name" column
# (column A) values, widens column A to fit, and leaves the active selection
# on the last edited row (A8), matching the author's final saved state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the dataset identifiers in column A (rows 2-8).
$ws.Range("A2").Value = "GSE65159_TSAI"
$ws.Range("A3").Value = "GDS4414_Muller"
$ws.Range("A4").Value = "GSE56772_BMS"
$ws.Range("A5").Value = "GSE57528_BMS"
$ws.Range("A6").Value = "GSE57583_BMS"
$ws.Range("A7").Value = "GSE9566_BARRES"
$ws.Range("A8").Value = "GSE31624_Pfizer"

# 2. Widen column A so the new (and old) names are fully visible.
#    40.1666... characters is stored by Excel as a column "width" of 41.
$ws.Columns.Item(1).ColumnWidth = 40.166666666666664

# 3. Leave the active cell / selection on A8 (last row touched).
$ws.Range("A8").Select()

# 4. Make sure the sheet prints in portrait orientation.
$ws.PageSetup.Orientation = 1
